# Apply updated airtoxics NATA data to the Means and Standard Deviations sheets.

$wb = $excel.ActiveWorkbook

# --- "Means" sheet ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 33
$wsMeans.Range("D9").Value = 67
$wsMeans.Range("E9").Value = 67
$wsMeans.Range("F9").Value = 65
$wsMeans.Range("G9").Value = 54

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.37
$wsMeans.Range("D10").Value = 0.42
$wsMeans.Range("E10").Value = 0.43
$wsMeans.Range("F10").Value = 0.43
$wsMeans.Range("G10").Value = 0.42

# --- "Standard Deviations" sheet ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million) SD
$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 9.5
$wsSD.Range("D9").Value = 8.2
$wsSD.Range("E9").Value = 7.1
$wsSD.Range("F9").Value = 6.6
$wsSD.Range("G9").Value = 8.3

# Row 10: Total Respiratory (hazard quotient) SD
$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.083
$wsSD.Range("G10").Value = 0.057
